# Update the grade-distribution counts on the "SECOND YEAR" sheet.
# The number of students previously recorded for A+ / A / A- are moved
# down into D+ / D / F, per the instructor's corrected submission.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("SECOND YEAR")

$ws.Range("E8").Value  = 0   # A+  (was 2)
$ws.Range("E9").Value  = 0   # A   (was 7)
$ws.Range("E10").Value = 0   # A-  (was 9)
$ws.Range("E17").Value = 9   # D+  (was 0)
$ws.Range("E18").Value = 8   # D   (was 1)
$ws.Range("E19").Value = 7   # F   (was 5)

$ws.Range("O16").Select() | Out-Null
